$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Lunes 27/05/2024"
$ws.Range("B3").Value = "Martes 28/05/2024"
$ws.Range("B4").Value = "Miércoles 29/05/2024"
$ws.Range("B5").Value = "Jueves 30/05/2024"
$ws.Range("B6").Value = "Viernes 31/05/2024"
$ws.Range("B7").Value = "Lunes 03/06/2024"
$ws.Range("B8").Value = "Martes 04/06/2024"
$ws.Range("B9").Value = "Miércoles 05/06/2024"
$ws.Range("B10").Value = "Jueves 06/06/2024"
$ws.Range("B11").Value = "Viernes 07/06/2024"
